# Generate Report for Handback
# Refresh the handoff/handback timestamps for the "4db86924-79a8-4086-a3cc-834dd8a0519b" entry.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 4db86924... row (row 3)
$wsOverview.Range("G3").Value = "2016-09-03 00:49:06"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the 4db86924... row (row 3)
$wsZhCn.Range("H3").Value = "2016-09-03 00:48:56"
$wsZhCn.Range("K3").Value = "2016-09-03 00:49:28"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the 4db86924... row (row 3)
$wsDeDe.Range("H3").Value = "2016-09-03 00:49:06"
$wsDeDe.Range("K3").Value = "2016-09-03 00:49:35"
